# Applies the "Add latest output file" edit to the simulation-output
# document: merges a stray blank paragraph into the following one,
# refreshes a couple of numeric/text values, and renames every
# "IM-III" label to "IM" (re-padding the surrounding whitespace so the
# fixed-width columns still line up), plus a few Normal-style tweaks.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Merge the empty first paragraph into the paragraph that follows
#    it (deletes the paragraph mark that separates them - same as
#    placing the cursor at the end of the empty paragraph and
#    pressing Delete).
# ---------------------------------------------------------------
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------
# 2) Plain text fixes / relabeling
# ---------------------------------------------------------------

$d.Content.Find.Execute(
    "/scratch/borgqvist/Dropbox/Work/Projects/symmetry_based_model_selsection_carcinogenesis/Code/symmetry_toolbox.py:185: RuntimeWarning: overflow encountered in exp",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/scratch/borgqvist/Dropbox/Work/Projects/symmetry_based_model_selsection_carcinogenesis/Code/symmetry_toolbox.py:180: RuntimeWarning: overflow encountered in exp",
    2)

$d.Content.Find.Execute(
    "The IM-III myeloma:     epsilon_IM_III_myeloma  =       0.767368329872",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The IM myeloma: epsilon_IM_myeloma      =       0.767368329872",
    2)

$d.Content.Find.Execute(
    "The IM-III colon:       epsilon_IM_III_colon    =       0.650429833712",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The IM colon:   epsilon_IM_colon        =       0.650429833712",
    2)

$d.Content.Find.Execute(
    "The IM-III CML: epsilon_IM_III_CML      =       0.554786599426",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The IM CML:     epsilon_IM_CML  =       0.554786599426",
    2)

$d.Content.Find.Execute(
    "IM-III   myeloma:        epsilon_scale  =       0.767368329872",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IM       myeloma:        epsilon_scale  =       0.767368329872",
    2)

$d.Content.Find.Execute(
    "IM-III   colon:  epsilon_scale  =       0.650429833712",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IM       colon:  epsilon_scale  =       0.650429833712",
    2)

$d.Content.Find.Execute(
    "IM-III   CML:    epsilon_scale  =       0.554786599426",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IM       CML:    epsilon_scale  =       0.554786599426",
    2)

$d.Content.Find.Execute(
    "IM-III epsilon scale:   0.650",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IM epsilon scale:       0.650",
    2)

$d.Content.Find.Execute(
    "IM-III, 0.3836842",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IM,     0.3836842",
    2)

# Standalone "IM-III" run (must run after the more specific matches
# above so it does not eat into them).
$d.Content.Find.Execute(
    "IM-III",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "IM",
    2)

# ---------------------------------------------------------------
# 3) Normal style tweaks (paragraph spacing / justification /
#    hyphenation suppression)
# ---------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0
Write-Host "style updated"
